# Apply the changes described by the diff to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5: update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"

# C8: Total Billed Amount goes from 0 to 3828.4
$ws.Range("C8").Value = 3828.4

# G10: Scope ID # value is cleared out
$ws.Range("G10").Value = ""

# H16:H23 - Pricing column for each line item goes from 0 to 478.55
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55
$ws.Range("H18").Value = 478.55
$ws.Range("H19").Value = 478.55
$ws.Range("H20").Value = 478.55
$ws.Range("H21").Value = 478.55
$ws.Range("H22").Value = 478.55
$ws.Range("H23").Value = 478.55

# H24: TOTAL row, sum of the pricing column
$ws.Range("H24").Value = 3828.400000000001
